$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Update the "as of" timestamp in the header row (A1) of every sheet
#    from 11:25 to 15:16 on 20/04/2021.
# ---------------------------------------------------------------------------
$wsMD410 = $wb.Worksheets.Item("MD410 Attendance")
$wsMD410.Range("A1").Value = "MD410 Registrees as of 20/04/2021 15:16"

$ws410E = $wb.Worksheets.Item("410E Attendance")
$ws410E.Range("A1").Value = "410E Registrees as of 20/04/2021 15:16"

$ws410W = $wb.Worksheets.Item("410W Attendance")
$ws410W.Range("A1").Value = "410W Registrees as of 20/04/2021 15:16"

$ws410EVoting = $wb.Worksheets.Item("410E Voting")
$ws410EVoting.Range("A1").Value = "410E Voting details as of 20/04/2021 15:16"

$ws410WVoting = $wb.Worksheets.Item("410W Voting")
$ws410WVoting.Range("A1").Value = "410W Voting details as of 20/04/2021 15:16"

# ---------------------------------------------------------------------------
# 2. "410E Attendance": remove the registree "Zucker, Leonie" (row 110).
#    Everything below shifts up by one row.
# ---------------------------------------------------------------------------
$ws410E.Rows("110").Delete()

# Refresh the trailing summary rows (now one row higher than before).
$ws410E.Cells.Item(125, 1).Value = "Number of attendees: 122"
$ws410E.Cells.Item(126, 1).Value = "Number of voters: 59"

# ---------------------------------------------------------------------------
# 3. "410W Attendance": add a new registree "Fouche, Jean" (Worcester club)
#    right before the existing "Fourie, Michele" row (old row 37), so the
#    new row becomes row 37 and everything below shifts down by one row.
# ---------------------------------------------------------------------------
$ws410W.Rows("37").Insert()

# Match the formatting used by the surrounding data rows.
$ws410W.Range("A38:E38").Copy()
$ws410W.Range("A37:E37").PasteSpecial(-4122)  # xlPasteFormats
$ws410W.Rows("37").RowHeight = 25

$ws410W.Cells.Item(37, 1).Value = "Fouché"
$ws410W.Cells.Item(37, 2).Value = "Jéan"
$ws410W.Cells.Item(37, 3).Value = "Worcester"
$ws410W.Cells.Item(37, 4).Value = "No"
$ws410W.Cells.Item(37, 5).Value = "Yes"

# Refresh the trailing summary rows (now one row lower than before).
$ws410W.Cells.Item(151, 1).Value = "Number of attendees: 148"
$ws410W.Cells.Item(152, 1).Value = "Number of voters: 55"

# ---------------------------------------------------------------------------
# 4. "410E Voting": the Milnerton club no longer has any voters (its sole
#    voter, Zucker Leonie, was removed above), so drop its row (old row 17).
# ---------------------------------------------------------------------------
$ws410EVoting.Rows("17").Delete()

# Refresh the trailing summary rows (now one row higher than before).
$ws410EVoting.Cells.Item(30, 1).Value = "Number of clubs: 27"
$ws410EVoting.Cells.Item(31, 1).Value = "Number of voters: 60"
